# Updated code till Inventory Dashboard
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # login
$ws2 = $wb.Worksheets.Item(2)   # incorrect_cred
$ws3 = $wb.Worksheets.Item(3)   # invalid_email
$ws4 = $wb.Worksheets.Item(4)   # empty_email
$ws5 = $wb.Worksheets.Item(5)   # empty_password

# ---------------------------------------------------------------------------
# Sheet1 (login) - add rows 4 & 5 with new credentials + mailto hyperlinks
# ---------------------------------------------------------------------------
$ws1.Range("A4").Value = "rohit13@latido.com.np"
$ws1.Range("B4").Value = "Frugal@123"
$ws1.Hyperlinks.Add($ws1.Range("B4"), "mailto:rohit13@latido.com.np")
$ws1.Range("C4").Value = "OVERVIEW"

$ws1.Range("A5").Value = "rohit13@latido.com.np"
$ws1.Range("B5").Value = "Frugal@123"
$ws1.Hyperlinks.Add($ws1.Range("B5"), "mailto:rohit13@latido.com.np")
$ws1.Range("C5").Value = "OVERVIEW"

# ---------------------------------------------------------------------------
# Sheet2 (incorrect_cred) - message text changed, two new rows added
# ---------------------------------------------------------------------------
$ws3.Range("C2").Copy()
$ws2.Range("C2").PasteSpecial(-4122)
$ws2.Range("C2").Value = "Please check username"

$ws3.Range("C2").Copy()
$ws2.Range("C3").PasteSpecial(-4122)
$ws2.Range("C3").Value = "Please check username"

$ws2.Range("A4").Value = "rohit13@latido.com.np"
$ws2.Range("B4").Value = "Test@123"
$ws2.Hyperlinks.Add($ws2.Range("B4"), "mailto:rohit13@latido.com.np")
$ws3.Range("C2").Copy()
$ws2.Range("C4").PasteSpecial(-4122)
$ws2.Range("C4").Value = "Please check username"

$ws2.Range("A5").Value = "rohit13@latido.com.np"
$ws2.Hyperlinks.Add($ws2.Range("B5"), "mailto:rohit13@latido.com.np", "", "", "Test@123")
$ws2.Range("B5").Value = "Test@124"
$ws3.Range("C2").Copy()
$ws2.Range("C5").PasteSpecial(-4122)
$ws2.Range("C5").Value = "Please check username"

# ---------------------------------------------------------------------------
# Sheet3 (invalid_email) - two new rows added
# ---------------------------------------------------------------------------
$ws3.Range("A4").Value = "frugal@latido.com"
$ws3.Hyperlinks.Add($ws3.Range("A4"), "mailto:frugal@latido.com")
$ws3.Range("B4").Value = "Test@1824"
$ws3.Range("C4").Value = "Please check username"

$ws3.Range("A5").Value = "f rugallatido.com.np"
$ws3.Range("B5").Value = "Test@1825"
$ws3.Range("C5").Value = "Please check username"

# ---------------------------------------------------------------------------
# Sheet4 (empty_email) - two new rows added (no A value)
# ---------------------------------------------------------------------------
$ws4.Range("B4").Value = "Test@3543"
$ws4.Range("C4").Value = "Please input your Email!"

$ws4.Range("B5").Value = "Test@5263"
$ws4.Range("C5").Value = "Please input your Email!"

# ---------------------------------------------------------------------------
# Sheet5 (empty_password) - two new rows added (no B value).
# A4 carries the "Hyperlink" cell style (copied from an actual hyperlink
# cell) but is not itself a real hyperlink, matching the target workbook.
# ---------------------------------------------------------------------------
$ws1.Range("B4").Copy()
$ws5.Range("A4").PasteSpecial(-4122)
$ws5.Range("A4").Value = "rohit13@latido.com.np"
$ws5.Range("C4").Value = "Please input your password!"

$ws5.Range("A5").Value = "rohit13@latido.com.np"
$ws5.Range("C5").Value = "Please input your password!"

# ---------------------------------------------------------------------------
# Selections per sheet
# ---------------------------------------------------------------------------
$ws1.Range("E7").Select()
$ws3.Range("A6").Select()
$ws4.Range("D13").Select()
$ws5.Range("D11").Select()

$ws2.Range("C10").Select()
$ws2.Activate()
